$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 10159.42
$ws.Range("B8").Value = 10092.81
$ws.Range("C8").Value = 305.24
$ws.Range("D8").Value = 307.24
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = 0.66
$ws.Range("G8").Value = 42613.766782407409
$ws.Range("H8").Value = $true

$ws.Range("A9").Value = 10177.709999999999
$ws.Range("B9").Value = 10159.42
$ws.Range("C9").Value = 307.68
$ws.Range("D9").Value = 308.24
$ws.Range("E9").Value = $false
$ws.Range("F9").Value = 0.18
$ws.Range("G9").Value = 42614.674826388888
$ws.Range("H9").Value = $true

$ws.Range("A10").Value = 10173.64
$ws.Range("B10").Value = 10177.709999999999
$ws.Range("C10").Value = 307.95999999999998
$ws.Range("D10").Value = 307.83
$ws.Range("E10").Value = $false
$ws.Range("F10").Value = -0.04
$ws.Range("G10").Value = 42615.751875000002
$ws.Range("H10").Value = $false
